$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034815561719928
$ws.Range("D2").Value = 1.035297447015192
$ws.Range("E2").Value = 1.042762738731573
$ws.Range("F2").Value = 1.050890842935073
$ws.Range("I2").Value = 1.034534133364597
$ws.Range("J2").Value = 1.039932740324115
$ws.Range("K2").Value = 1.038094455115816
$ws.Range("L2").Value = 1.045538509624168
$ws.Range("M2").Value = 1.053643863112506
$ws.Range("N2").Value = 1.017162448885601

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036723362601519
$ws.Range("D3").Value = 1.035885561459303
$ws.Range("E3").Value = 1.044490844565024
$ws.Range("F3").Value = 1.05280599178856
$ws.Range("I3").Value = 1.034769907342393
$ws.Range("J3").Value = 1.041478991770305
$ws.Range("K3").Value = 1.038492523902997
$ws.Range("L3").Value = 1.04707511352033
$ws.Range("M3").Value = 1.055368712701102
$ws.Range("N3").Value = 1.01767990087845

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037950091219948
$ws.Range("D4").Value = 1.036261271265468
$ws.Range("E4").Value = 1.045600844316189
$ws.Range("F4").Value = 1.054032842284199
$ws.Range("I4").Value = 1.03491649139569
$ws.Range("J4").Value = 1.042471626739734
$ws.Range("K4").Value = 1.038744575414643
$ws.Range("L4").Value = 1.048060769489872
$ws.Range("M4").Value = 1.056472043143198
$ws.Range("N4").Value = 1.018011922914198

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03846398425653
$ws.Range("D5").Value = 1.036418064424282
$ws.Range("E5").Value = 1.046065552928426
$ws.Range("F5").Value = 1.054545680688036
$ws.Range("I5").Value = 1.034976687745506
$ws.Range("J5").Value = 1.042887066782579
$ws.Range("K5").Value = 1.038849218412895
$ws.Range("L5").Value = 1.048473099481937
$ws.Range("M5").Value = 1.056932859515044
$ws.Range("N5").Value = 1.018150842421148

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038550163097499
$ws.Range("D6").Value = 1.036444323032604
$ws.Range("E6").Value = 1.046143466859169
$ws.Range("F6").Value = 1.054631617690058
$ws.Range("I6").Value = 1.034986711356661
$ws.Range("J6").Value = 1.042956712446053
$ws.Range("K6").Value = 1.038866711151093
$ws.Range("L6").Value = 1.048542212634044
$ws.Range("M6").Value = 1.057010056197283
$ws.Range("N6").Value = 1.018174129014341

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037956965008603
$ws.Range("D7").Value = 1.03626337087801
$ws.Range("E7").Value = 1.045607061342559
$ws.Range("F7").Value = 1.054039706321897
$ws.Range("I7").Value = 1.034917301347047
$ws.Range("J7").Value = 1.042477185158395
$ws.Range("K7").Value = 1.038745978839442
$ws.Range("L7").Value = 1.048066287043479
$ws.Range("M7").Value = 1.056478212433229
$ws.Range("N7").Value = 1.018013781755193

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035461934628642
$ws.Range("D8").Value = 1.035497207245904
$ws.Range("E8").Value = 1.043348475161851
$ws.Range("F8").Value = 1.051540660414111
$ws.Range("I8").Value = 1.034615054137837
$ws.Range("J8").Value = 1.040456953413713
$ws.Range("K8").Value = 1.038230130592062
$ws.Range("L8").Value = 1.046059614899318
$ws.Range("M8").Value = 1.054229446094204
$ws.Range("N8").Value = 1.017337910131062

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03100462721931
$ws.Range("D9").Value = 1.034109893198158
$ws.Range("E9").Value = 1.039304462225264
$ws.Range("F9").Value = 1.047040632927376
$ws.Range("I9").Value = 1.034036502637194
$ws.Range("J9").Value = 1.036835395170628
$ws.Range("K9").Value = 1.037278649100277
$ws.Range("L9").Value = 1.042456297969375
$ws.Range("M9").Value = 1.050167610431851
$ws.Range("N9").Value = 1.016125059080258

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027990208780578
$ws.Range("D10").Value = 1.033159712661439
$ws.Range("E10").Value = 1.036563496932742
$ws.Range("F10").Value = 1.043973514424503
$ws.Range("I10").Value = 1.033619634612117
$ws.Range("J10").Value = 1.034377842122566
$ws.Range("K10").Value = 1.036615490036467
$ws.Range("L10").Value = 1.040007083406405
$ws.Range("M10").Value = 1.047390864392556
$ws.Range("N10").Value = 1.01530120031428

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026674292436487
$ws.Range("D11").Value = 1.032742203364343
$ws.Range("E11").Value = 1.035365534749142
$ws.Range("F11").Value = 1.042628969256411
$ws.Range("I11").Value = 1.033431665293385
$ws.Range("J11").Value = 1.033303052808874
$ws.Range("K11").Value = 1.036321426279513
$ws.Range("L11").Value = 1.038934991967482
$ws.Range("M11").Value = 1.046171670837178
$ws.Range("N11").Value = 1.014940699230803

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026183860528799
$ws.Range("D12").Value = 1.032586202569772
$ws.Range("E12").Value = 1.034918851480594
$ws.Range("F12").Value = 1.042127024439101
$ws.Range("I12").Value = 1.033360717419134
$ws.Range("J12").Value = 1.032902192333577
$ws.Range("K12").Value = 1.036211152988595
$ws.Range("L12").Value = 1.03853499585217
$ws.Range("M12").Value = 1.04571623275531
$ws.Range("N12").Value = 1.014806215566526

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026289134819485
$ws.Range("D13").Value = 1.032619706978097
$ws.Range("E13").Value = 1.035014744419184
$ws.Range("F13").Value = 1.042234808161885
$ws.Range("I13").Value = 1.033375987118812
$ws.Range("J13").Value = 1.032988252873279
$ws.Range("K13").Value = 1.036234854378232
$ws.Range("L13").Value = 1.038620877227864
$ws.Range("M13").Value = 1.045814043041731
$ws.Range("N13").Value = 1.014835089097351

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02663378692227
$ws.Range("D14").Value = 1.032729327082095
$ws.Range("E14").Value = 1.035328646813903
$ws.Range("F14").Value = 1.042587530047384
$ws.Range("I14").Value = 1.0334258237639
$ws.Range("J14").Value = 1.033269951137172
$ws.Range("K14").Value = 1.036312332410325
$ws.Range("L14").Value = 1.038901964557943
$ws.Range("M14").Value = 1.046134076982932
$ws.Range("N14").Value = 1.014929594620322

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0268459196174
$ws.Range("D15").Value = 1.032796745631959
$ws.Range("E15").Value = 1.035521825111378
$ws.Range("F15").Value = 1.042804518150852
$ws.Range("I15").Value = 1.033456380132481
$ws.Range("J15").Value = 1.03344329688708
$ws.Range("K15").Value = 1.036359930517938
$ws.Range("L15").Value = 1.039074915618372
$ws.Range("M15").Value = 1.046330917839426
$ws.Range("N15").Value = 1.014987745704276

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028077313851869
$ws.Range("D16").Value = 1.033187292837283
$ws.Range("E16").Value = 1.036642764510735
$ws.Range("F16").Value = 1.04406239658902
$ws.Range("I16").Value = 1.033631951739803
$ws.Range("J16").Value = 1.034448944911817
$ws.Range("K16").Value = 1.036634859913527
$ws.Range("L16").Value = 1.040077987864517
$ws.Range("M16").Value = 1.04747141958965
$ws.Range("N16").Value = 1.01532504525716

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02884685385885
$ws.Range("D17").Value = 1.033430641869079
$ws.Range("E17").Value = 1.037342900165164
$ws.Range("F17").Value = 1.044846989934177
$ws.Range("I17").Value = 1.033740080726752
$ws.Range("J17").Value = 1.035076884672134
$ws.Range("K17").Value = 1.036805460971638
$ws.Range("L17").Value = 1.040704067120565
$ws.Range("M17").Value = 1.048182285505323
$ws.Range("N17").Value = 1.015535608130805

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029294688504554
$ws.Range("D18").Value = 1.033571997577309
$ws.Range("E18").Value = 1.037750208512862
$ws.Range("F18").Value = 1.045303044855638
$ws.Range("I18").Value = 1.033802430956585
$ws.Range("J18").Value = 1.035442125939351
$ws.Range("K18").Value = 1.036904303290076
$ws.Range("L18").Value = 1.041068135166205
$ws.Range("M18").Value = 1.048595299449585
$ws.Range("N18").Value = 1.015658063466718

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029447215919971
$ws.Range("D19").Value = 1.033620097073367
$ws.Range("E19").Value = 1.037888909961193
$ws.Range("F19").Value = 1.045458280350888
$ws.Range("I19").Value = 1.033823568910652
$ws.Range("J19").Value = 1.035566490907123
$ws.Range("K19").Value = 1.036937893125338
$ws.Range("L19").Value = 1.041192085183621
$ws.Range("M19").Value = 1.048735852718709
$ws.Range("N19").Value = 1.015699756452219

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028764395817984
$ws.Range("D20").Value = 1.033404593445423
$ws.Range("E20").Value = 1.037267893073191
$ws.Range("F20").Value = 1.044762974787347
$ws.Range("I20").Value = 1.033728553985105
$ws.Range("J20").Value = 1.035009618990648
$ws.Range("K20").Value = 1.036787226066878
$ws.Range("L20").Value = 1.040637010162302
$ws.Range("M20").Value = 1.048106184421785
$ws.Range("N20").Value = 1.01551305429958

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026532341134023
$ws.Range("D21").Value = 1.032697072124528
$ws.Range("E21").Value = 1.035236257810582
$ws.Range("F21").Value = 1.042483732221965
$ws.Range("I21").Value = 1.033411179298496
$ws.Range("J21").Value = 1.033187043467616
$ws.Range("K21").Value = 1.036289545960525
$ws.Range("L21").Value = 1.038819240601064
$ws.Range("M21").Value = 1.046039906421545
$ws.Range("N21").Value = 1.014901781138718

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02511943831338
$ws.Range("D22").Value = 1.032246903301043
$ws.Range("E22").Value = 1.033948993409828
$ws.Range("F22").Value = 1.041036072198406
$ws.Range("I22").Value = 1.033205105180018
$ws.Range("J22").Value = 1.032031634905193
$ws.Range("K22").Value = 1.035970585592437
$ws.Range("L22").Value = 1.037666056400572
$ws.Range("M22").Value = 1.044725832202587
$ws.Range("N22").Value = 1.014514101986201

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025869361515181
$ws.Range("D23").Value = 1.032486053094797
$ws.Range("E23").Value = 1.03463234759764
$ws.Range("F23").Value = 1.041804905309641
$ws.Range("I23").Value = 1.033314970037931
$ws.Range("J23").Value = 1.032645049843316
$ws.Range("K23").Value = 1.036140248204193
$ws.Range("L23").Value = 1.038278367964468
$ws.Range("M23").Value = 1.045423877684343
$ws.Range("N23").Value = 1.01471993940614

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028801658222824
$ws.Range("D24").Value = 1.033416365416554
$ws.Range("E24").Value = 1.037301788850802
$ws.Range("F24").Value = 1.044800942511711
$ws.Range("I24").Value = 1.033733764646852
$ws.Range("J24").Value = 1.035040016623571
$ws.Range("K24").Value = 1.036795467694514
$ws.Range("L24").Value = 1.040667313754143
$ws.Range("M24").Value = 1.048140576236106
$ws.Range("N24").Value = 1.015523246524483

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032164352601877
$ws.Range("D25").Value = 1.03447298592578
$ws.Range("E25").Value = 1.040357719978881
$ws.Range("F25").Value = 1.048215651267788
$ws.Range("I25").Value = 1.034191544136822
$ws.Range("J25").Value = 1.037779135333547
$ws.Range("K25").Value = 1.037529691485739
$ws.Range("L25").Value = 1.043395994198505
$ws.Range("M25").Value = 1.051229666170977
$ws.Range("N25").Value = 1.016441261573967
